$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 data: switch the test-data row to Preprod / password / Netherlands ---
$ws.Range("A2").Value = "Preprod"
$ws.Range("C2").Value = "password"
$ws.Range("D2").Value = "Netherlands"

# Give D2 the same "wrap + vertical-center" formatting already used by the
# other Country cells (D9/D11) - copy format only, value was already set above.
$ws.Range("D9").Copy()
$ws.Range("D2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 2 grew taller once the Country column started wrapping
$ws.Rows.Item(2).RowHeight = 28.8

# New Country column (D) gets its own width
$ws.Columns.Item(4).ColumnWidth = 13.25

# Scrolled / selected a different cell afterwards
$ws.Range("F23").Select()
